$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1. Update the Temperature row (row 7, before the new row is inserted) with
#    the new NOAA/netcdf sourced content.
$ws.Range("D7").Value = "Using 0.5 degree by 0.5 degree temperature projections (50 km x 50 kmish cells).`nTemperature will vary monthly and interact with local-level variables such that the`nthe local-level variable's effects on movement can vary with temperature. "

$ws.Range("B7").Value = "NOAA data`nhttps://www.esrl.noaa.gov/psd/data/gridded/data.ghcncams.html"
$ws.Range("B7").WrapText = $true

$ws.Range("E7").Value = "Time-varying, but constant at the `nstudy level."
$ws.Range("E7").WrapText = $true

# 2. Insert a brand-new "Elevation" row right after "Distance to water" (row 6),
#    pushing Temperature and everything below it down by one row.
$ws.Rows.Item(7).Insert()
$ws.Range("A7").Value = "Elevation"

# 3. Header "Level" -> "Level/Dynamics"
$ws.Range("E1").Value = "Level/Dynamics"

# 4. Plant productivity row: "Cell-level" -> "Cell-level/Varies monthly"
$ws.Range("E3").Value = "Cell-level/Varies monthly"

# 5. Finish off the new Elevation row's remaining cell + row height.
$ws.Range("E7").Value = "Cell-level"
$ws.Rows.Item(7).RowHeight = 58

# 6. Update the view: scroll back to the top and move the selection to E10.
$ws.Application.ActiveWindow.ScrollRow = 1
$ws.Application.ActiveWindow.ScrollColumn = 1
$ws.Range("E10").Select()
